# Teacher availability -> preferred timeslots refactor (Teachers sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teachers")

# Header: availableTimeslots -> preferredTimeslots
$ws.Cells.Item(1, 3).Value = "preferredTimeslots"

# Row 4 (I. Jones): "2, 5" -> "TUESDAY/10:00-16:00|WEDNESDAY/09:00-15:00"
# and drop the custom number-format style that cell used to carry.
$ws.Cells.Item(4, 3).ClearFormats()
$ws.Cells.Item(4, 3).Value = "TUESDAY/10:00-16:00|WEDNESDAY/09:00-15:00"

# Row 2 (A. Turing): "1,2,4,5" -> "MONDAY/08:00-12:00"
$ws.Cells.Item(2, 3).Value = "MONDAY/08:00-12:00"

# Widen column C to fit the longer preferred-timeslot strings.
$ws.Columns(3).ColumnWidth = 40.1666666666667

# Move the remembered selection to H7 (minor UI tweak from the commit).
$ws.Range("H7").Select()
